$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 13 values (C13, D13)
$ws.Range("C13").Value = 1064
$ws.Range("D13").Value = 39151

# Add new row 14: Sergipe / 01/01/2022 / 1104 / 40163
$ws.Range("A14").Value = "Sergipe"

# B14 must stay a plain text string like the other "date" cells in column B
# (they are inline/shared strings, not real Excel dates), so force text
# formatting before the assignment, then reset the style back to Normal so
# no stray style index is left on the cell.
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "01/01/2022"
$ws.Range("B14").Style = "Normal"

$ws.Range("C14").Value = 1104
$ws.Range("D14").Value = 40163
